$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule "R30" (row 10): the "Integer min" / "From" value in C10 changes from 18 to 1
$ws.Range("C10").Value = 1
